$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row layout: A=index(unchanged), B=random id(changed), C=image(unchanged text->reordered value),
# D=word(unchanged text->reordered value), E=category(recomputed to match image)

$ws.Cells.Item(2, 2).Value = 103
$ws.Cells.Item(2, 3).Value = "flower/flower027.jpg"
$ws.Cells.Item(2, 4).Value = "biegen"
$ws.Cells.Item(2, 5).Value = "flower"

$ws.Cells.Item(3, 2).Value = 43
$ws.Cells.Item(3, 3).Value = "flower/flower005.jpg"
$ws.Cells.Item(3, 4).Value = "raten"
$ws.Cells.Item(3, 5).Value = "flower"

$ws.Cells.Item(4, 2).Value = 108
$ws.Cells.Item(4, 3).Value = "dog/dog030.jpg"
$ws.Cells.Item(4, 4).Value = "fliehen"
$ws.Cells.Item(4, 5).Value = "dog"

$ws.Cells.Item(5, 2).Value = 109
$ws.Cells.Item(5, 3).Value = "flower/flower004.jpg"
$ws.Cells.Item(5, 4).Value = "fliegen"
$ws.Cells.Item(5, 5).Value = "flower"

$ws.Cells.Item(6, 2).Value = 110
$ws.Cells.Item(6, 3).Value = "dog/dog020.jpg"
$ws.Cells.Item(6, 4).Value = "saufen"
$ws.Cells.Item(6, 5).Value = "dog"

$ws.Cells.Item(7, 2).Value = 72
$ws.Cells.Item(7, 3).Value = "dog/dog004.jpg"
$ws.Cells.Item(7, 4).Value = "enden"
$ws.Cells.Item(7, 5).Value = "dog"

$ws.Cells.Item(8, 2).Value = 86
$ws.Cells.Item(8, 3).Value = "flower/flower029.jpg"
$ws.Cells.Item(8, 4).Value = "antun"
$ws.Cells.Item(8, 5).Value = "flower"

$ws.Cells.Item(9, 2).Value = 25
$ws.Cells.Item(9, 3).Value = "dog/dog028.jpg"
$ws.Cells.Item(9, 4).Value = "langen"
$ws.Cells.Item(9, 5).Value = "dog"

$ws.Cells.Item(10, 2).Value = 6
$ws.Cells.Item(10, 3).Value = "dog/dog013.jpg"
$ws.Cells.Item(10, 4).Value = "wenden"
$ws.Cells.Item(10, 5).Value = "dog"

$ws.Cells.Item(11, 2).Value = 122
$ws.Cells.Item(11, 3).Value = "dog/dog015.jpg"
$ws.Cells.Item(11, 4).Value = "rücken"
$ws.Cells.Item(11, 5).Value = "dog"

$ws.Cells.Item(12, 2).Value = 90
$ws.Cells.Item(12, 3).Value = "dog/dog000.jpg"
$ws.Cells.Item(12, 4).Value = "schmecken"
$ws.Cells.Item(12, 5).Value = "dog"

$ws.Cells.Item(13, 2).Value = 68
$ws.Cells.Item(13, 3).Value = "flower/flower009.jpg"
$ws.Cells.Item(13, 4).Value = "stärken"
$ws.Cells.Item(13, 5).Value = "flower"

$ws.Cells.Item(14, 2).Value = 62
$ws.Cells.Item(14, 3).Value = "flower/flower006.jpg"
$ws.Cells.Item(14, 4).Value = "lehnen"
$ws.Cells.Item(14, 5).Value = "flower"

$ws.Cells.Item(15, 2).Value = 18
$ws.Cells.Item(15, 3).Value = "flower/flower031.jpg"
$ws.Cells.Item(15, 4).Value = "opfern"
$ws.Cells.Item(15, 5).Value = "flower"

$ws.Cells.Item(16, 2).Value = 84
$ws.Cells.Item(16, 3).Value = "dog/dog029.jpg"
$ws.Cells.Item(16, 4).Value = "mieten"
$ws.Cells.Item(16, 5).Value = "dog"

$ws.Cells.Item(17, 2).Value = 48
$ws.Cells.Item(17, 3).Value = "flower/flower025.jpg"
$ws.Cells.Item(17, 4).Value = "formen"
$ws.Cells.Item(17, 5).Value = "flower"

$ws.Cells.Item(18, 2).Value = 92
$ws.Cells.Item(18, 3).Value = "dog/dog007.jpg"
$ws.Cells.Item(18, 4).Value = "füttern"
$ws.Cells.Item(18, 5).Value = "dog"

$ws.Cells.Item(19, 2).Value = 4
$ws.Cells.Item(19, 3).Value = "flower/flower008.jpg"
$ws.Cells.Item(19, 4).Value = "pflegen"
$ws.Cells.Item(19, 5).Value = "flower"

$ws.Cells.Item(20, 2).Value = 12
$ws.Cells.Item(20, 3).Value = "dog/dog017.jpg"
$ws.Cells.Item(20, 4).Value = "liefern"
$ws.Cells.Item(20, 5).Value = "dog"

$ws.Cells.Item(21, 2).Value = 40
$ws.Cells.Item(21, 3).Value = "flower/flower000.jpg"
$ws.Cells.Item(21, 4).Value = "tauschen"
$ws.Cells.Item(21, 5).Value = "flower"

$ws.Cells.Item(22, 2).Value = 53
$ws.Cells.Item(22, 3).Value = "flower/flower015.jpg"
$ws.Cells.Item(22, 4).Value = "backen"
$ws.Cells.Item(22, 5).Value = "flower"

$ws.Cells.Item(23, 2).Value = 91
$ws.Cells.Item(23, 3).Value = "flower/flower001.jpg"
$ws.Cells.Item(23, 4).Value = "segeln"
$ws.Cells.Item(23, 5).Value = "flower"

$ws.Cells.Item(24, 2).Value = 28
$ws.Cells.Item(24, 3).Value = "dog/dog001.jpg"
$ws.Cells.Item(24, 4).Value = "runden"
$ws.Cells.Item(24, 5).Value = "dog"

$ws.Cells.Item(25, 2).Value = 47
$ws.Cells.Item(25, 3).Value = "dog/dog023.jpg"
$ws.Cells.Item(25, 4).Value = "stechen"
$ws.Cells.Item(25, 5).Value = "dog"

$ws.Cells.Item(26, 2).Value = 14
$ws.Cells.Item(26, 3).Value = "flower/flower016.jpg"
$ws.Cells.Item(26, 4).Value = "schicken"
$ws.Cells.Item(26, 5).Value = "flower"

$ws.Cells.Item(27, 2).Value = 98
$ws.Cells.Item(27, 3).Value = "dog/dog006.jpg"
$ws.Cells.Item(27, 4).Value = "fesseln"
$ws.Cells.Item(27, 5).Value = "dog"

$ws.Cells.Item(28, 2).Value = 13
$ws.Cells.Item(28, 3).Value = "dog/dog002.jpg"
$ws.Cells.Item(28, 4).Value = "drohen"
$ws.Cells.Item(28, 5).Value = "dog"

$ws.Cells.Item(29, 2).Value = 123
$ws.Cells.Item(29, 3).Value = "flower/flower007.jpg"
$ws.Cells.Item(29, 4).Value = "tagen"
$ws.Cells.Item(29, 5).Value = "flower"

$ws.Cells.Item(30, 2).Value = 42
$ws.Cells.Item(30, 3).Value = "dog/dog010.jpg"
$ws.Cells.Item(30, 4).Value = "nehmen"
$ws.Cells.Item(30, 5).Value = "dog"

$ws.Cells.Item(31, 2).Value = 95
$ws.Cells.Item(31, 3).Value = "flower/flower026.jpg"
$ws.Cells.Item(31, 4).Value = "jubeln"
$ws.Cells.Item(31, 5).Value = "flower"

$ws.Cells.Item(32, 2).Value = 55
$ws.Cells.Item(32, 3).Value = "flower/flower003.jpg"
$ws.Cells.Item(32, 4).Value = "strahlen"
$ws.Cells.Item(32, 5).Value = "flower"

$ws.Cells.Item(33, 2).Value = 66
$ws.Cells.Item(33, 3).Value = "dog/dog014.jpg"
$ws.Cells.Item(33, 4).Value = "sondern"
$ws.Cells.Item(33, 5).Value = "dog"
